$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "additionId"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "nameAr"
$ws.Range("D1").Value = "phone"

# --- Data row (row 2) ---
$ws.Range("A2").Value = 2001
$ws.Range("B2").Value = "MADI INTERNATIONNAL"
$ws.Range("C2").Value = "MADI INTERNATIONNAL"
$ws.Range("D2").Value = 5005001

# --- Column widths (mirrors the bestFit widths baked into the authored file) ---
$ws.Columns("A").ColumnWidth = 10.140625
$ws.Columns("B:C").ColumnWidth = 22.5703125
$ws.Columns("N").ColumnWidth = 21.42578125

# --- Final selection left on A2, as in the authored workbook ---
[void]$ws.Range("A2").Select()
